$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.849.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +7.77%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.864.17'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +14.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -1.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '426.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +12.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +10.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.858.62'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +15.15%  '
$ws.Range('E8').Value = '  +7.21%  '
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.724'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.42%  '
$ws.Range('E11').Value = '  +20.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000344'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +26.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.96'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +10.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.494.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +13.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '10.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +15.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +34.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.892.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +14.23%  '
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.044.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.25%  '
$ws.Range('E21').Value = '  +10.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '412.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +9.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.87'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +14.87%  '
$ws.Range('E24').Value = '  +8.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.59'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +19.16%  '
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +18.76%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.03%  '
$ws.Range('E29').Value = '  +3.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +40.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '730.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +14.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.60'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +19.83%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.122'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +15.64%  '
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.70'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '39.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +11.77%  '
$ws.Range('E37').Value = '  +4.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.91'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.89%  '
$ws.Range('E39').Value = '  +37.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0759'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +35.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0459'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.95%  '
$ws.Range('E42').Value = '  +12.71%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  +9.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.135'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +14.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.313'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +18.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '142.27'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('E49').Value = '  +9.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.81'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.95%  '
$ws.Range('E51').Value = '  +5.98%  '
